# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet Hoja1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.88 = 10565.75 pesos`n✅ 10565.75 pesos = 2.86 = 957.52 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the tasas sheet values (N10/O10, N12/O12) ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 347.822
$wsTasas.Range("O10").Value = 3675
$wsTasas.Range("N12").Value = 3700
$wsTasas.Range("O12").Value = 335.311
